$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" block
# ("LOM3013: Ciência dos Materiais (Requisito)") and the paragraph that
# contains the site's footer/copyright notice. Everything between the end
# of the former and the end of the latter (an empty paragraph, a
# page-break paragraph, and the copyright paragraph itself) is removed,
# leaving the trailing empty paragraph + page-break paragraph intact.

$count = $d.Paragraphs.Count

$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOM3013: Ciência dos Materiais (Requisito)*") {
        $startIdx = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -gt $startIdx) {
    $delStart = $d.Paragraphs.Item($startIdx).Range.End
    $delEnd = $d.Paragraphs.Item($endIdx).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
